# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets carry identical data in this workbook, so the same set of
# row -> new value updates is applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2184
    5  = 12819
    6  = 66
    8  = 507
    9  = 471
    10 = 1161
    11 = 964
    12 = 13694
    13 = 14137
    18 = 25
    22 = 1073
    25 = 933
    26 = 5234
    28 = 278
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
